$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in this sheet are stored as text (e.g. "287.31", "1.37%") even
# though they look numeric. A plain `Range.Value = "287.46"` assignment
# would have Excel auto-convert the literal into a Number cell, which
# changes the stored type. To keep the cell as Text (matching the
# original authoring) we briefly force a Text number format before the
# assignment, then restore the cell's style to "Normal" afterwards so we
# don't leave a stray formatting change behind (the value itself stays
# text once stored, switching the display format back to General does
# not re-parse it).
function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "287.46"
Set-TextValue "E2" "1.42%"
Set-TextValue "D3" "29.56"
Set-TextValue "E3" "3.98%"
Set-TextValue "D4" "5.134"
Set-TextValue "E4" "1.54%"
Set-TextValue "D5" "0.06707"
Set-TextValue "E5" "3.21%"
Set-TextValue "D6" "7.342"
Set-TextValue "E6" "1.64%"
Set-TextValue "D7" "3.397"
Set-TextValue "E7" "1.07%"
Set-TextValue "E8" "-2.78%"
Set-TextValue "D9" "0.9204"
Set-TextValue "E9" "0.29%"
Set-TextValue "D10" "0.1590"
Set-TextValue "E10" "3.58%"
Set-TextValue "D11" "0.06839"
Set-TextValue "E11" "3.53%"
Set-TextValue "D12" "0.07704"
Set-TextValue "E12" "1.21%"
Set-TextValue "D13" "0.02938"
Set-TextValue "E13" "5.19%"
Set-TextValue "D14" "0.08986"
Set-TextValue "E14" "0.20%"
Set-TextValue "D15" "0.001589"
Set-TextValue "E15" "-0.08%"
Set-TextValue "D16" "0.04495"
Set-TextValue "E16" "1.46%"
Set-TextValue "D17" "0.0006449"
Set-TextValue "E17" "0.87%"
Set-TextValue "D18" "0.006251"
Set-TextValue "E18" "2.70%"
Set-TextValue "E19" "-0.03%"
Set-TextValue "E20" "-0.42%"
Set-TextValue "D21" "0.3216"
Set-TextValue "E21" "1.09%"
Set-TextValue "E22" "-2.84%"
Set-TextValue "D23" "4.073"
Set-TextValue "E23" "1.51%"
Set-TextValue "E24" "2.45%"
Set-TextValue "D25" "0.001193"
Set-TextValue "E25" "0.62%"
Set-TextValue "D26" "0.004122"
Set-TextValue "E26" "-7.56%"
Set-TextValue "D27" "0.0001198"
Set-TextValue "E27" "-0.18%"
Set-TextValue "D28" "0.0001616"
Set-TextValue "E28" "-0.15%"
Set-TextValue "D40" "0.04270"
Set-TextValue "E40" "3.76%"
Set-TextValue "D41" "0.006735"
Set-TextValue "D42" "0.1241"
Set-TextValue "E42" "0.61%"
Set-TextValue "D43" "0.002206"
Set-TextValue "E43" "5.55%"
Set-TextValue "D44" "0.01210"
Set-TextValue "E44" "5.02%"
Set-TextValue "D45" "0.00005684"
Set-TextValue "E45" "5.36%"
Set-TextValue "D46" "1.968"
Set-TextValue "E46" "-3.59%"
Set-TextValue "E47" "-29.47%"
